# Handle Disconnects; Server Complete?
#
# Adds a new "Disconnect" row (row 44) to the bottom of the message table
# on Sheet1: Category column is empty (it belongs to the "New Game"
# section started at row 41), Client Msg (C) = "Disconnect",
# Server Msg (D) = "disconn", tag/info (E) left blank, Comments (F) =
# "End game" -- mirroring the layout of the existing rows 42/43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C44").Value = "Disconnect"
$ws.Range("D44").Value = "disconn"
$ws.Range("F44").Value = "End game"

# Move the view/selection to reflect scrolling down to the newly added row,
# the way Excel would leave things after typing the new row in at the
# bottom of the sheet.
[void]$ws.Range("C45").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
